$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 (Logistic Regression) values
$ws.Range("B2").Value = 0.7548209366391185
$ws.Range("C2").Value = 0.7596885878811692
$ws.Range("D2").Value = 0.7548209366391185
$ws.Range("E2").Value = 0.7439024601770136

# Update row 3: change model name from Lasso to LightGBM, and update values to old LightGBM row values
$ws.Range("A3").Value = "LightGBM"
$ws.Range("B3").Value = 0.7870523415977961
$ws.Range("C3").Value = 0.7862886492881731
$ws.Range("D3").Value = 0.7870523415977961
$ws.Range("E3").Value = 0.7832676123203658

# Delete rows 4 through 8 (Support Vector Classifier, CART, Random Forest, old LightGBM, XGBoost)
$ws.Rows("4:8").Delete()
